# Adult TOD norms update ("adult TOD norms through pan").
#
# Each of the 6 age-band sheets holds a two-column raw-score -> standard-score
# lookup table (A: raw, B: ss) running from row 2 to row 21 (raw 1-20).
# This adds a raw=0 entry at the top and renumbers/rescales the rest of the
# table down one row, so raw now runs 0-20 (rows 2-22), with new ss values.

$wb = $excel.ActiveWorkbook

# New standard-score (column B) values for raw = 0..20, per sheet (in
# worksheet order: 18.0-23.11, 24.0-39.11, 40.0-49.11, 50.0-59.11,
# 60.0-69.11, 70.0-89.11).
$ssBySheet = @(
    @(46,52,57,63,69,75,81,86,92,98,104,110,115,121,127,130,130,130,130,130,130),
    @(48,54,60,66,71,77,83,89,95,101,106,112,118,124,129,130,130,130,130,130,130),
    @(50,56,62,68,73,79,85,91,96,102,108,113,119,124,129,130,130,130,130,130,130),
    @(51,57,63,69,74,80,86,91,97,102,107,113,118,123,128,130,130,130,130,130,130),
    @(52,58,64,70,75,81,86,92,97,102,107,112,117,122,127,130,130,130,130,130,130),
    @(54,60,66,72,78,83,89,94,99,104,109,114,119,124,128,130,130,130,130,130,130)
)

for ($sheetIdx = 1; $sheetIdx -le 6; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $ssValues = $ssBySheet[$sheetIdx - 1]

    for ($i = 0; $i -lt $ssValues.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $i
        $ws.Cells.Item($row, 2).Value = $ssValues[$i]
    }
}
